# Add "better support for importing capital calls & distributions":
# two new columns - "Generate Payments" (F) and "Payments Paid" (G) -
# with Yes/No flags per distribution row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("F1").Value = "Generate Payments"
$ws.Range("G1").Value = "Payments Paid"

# Row 2 (Aug distribution)
$ws.Range("F2").Value = "Yes"
$ws.Range("G2").Value = "Yes"

# Row 3 (Sept distribution)
$ws.Range("F3").Value = "Yes"
$ws.Range("G3").Value = "No"

# Row 4 (Nov distribution)
$ws.Range("F4").Value = "No"
$ws.Range("G4").Value = "No"

# Size the new columns to fit their content, like Excel's best-fit does
# automatically when data is typed into a previously empty column.
$ws.Columns.Item(6).ColumnWidth = 16.1
$ws.Columns.Item(7).ColumnWidth = 12.1

# Leave the selection where the user's cursor would land after filling
# in the new table columns.
$ws.Range("F5").Select()
